# Insert a new daily price record at row 571 (pushing the existing
# rows 571-638 down to 572-639), for:
#   Hortaliza, Vega Central Mapocho de Santiago - Ciboulette
#
# The new row copies the surrounding record's fixed attributes and only
# carries a new date (Fecha) and Volumen value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 571..638 down by one, leaving a blank row 571 to fill in.
$ws.Rows.Item(571).Insert()

$ws.Cells.Item(571, 1).Value = 9
$ws.Cells.Item(571, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(571, 3).Value = "Metropolitana"
$ws.Cells.Item(571, 4).Value = 45212
$ws.Cells.Item(571, 5).Value = 13
$ws.Cells.Item(571, 6).Value = 100112039
$ws.Cells.Item(571, 7).Value = "Ciboulette"
$ws.Cells.Item(571, 8).Value = "Sin especificar"
$ws.Cells.Item(571, 9).Value = "Primera"
$ws.Cells.Item(571, 10).Value = 430
$ws.Cells.Item(571, 11).Value = 1000
$ws.Cells.Item(571, 12).Value = 1200
$ws.Cells.Item(571, 13).Value = 1100
$ws.Cells.Item(571, 14).Value = "`$/docena de atados"
$ws.Cells.Item(571, 15).Value = "Región Metropolitana"
$ws.Cells.Item(571, 16).Value = 367
$ws.Cells.Item(571, 17).Value = 3
$ws.Cells.Item(571, 18).Value = "Hortaliza"
